$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = '[(-2.0, 0.0), (-2.0, -0.13), (-13.452, -0.2724750527055538), (-13.452, -0.4024750527055538), (-2.0, -0.13), (-2.0, -0.25), (-13.452, -0.4024750527055538), (-13.452, -0.5224750527055537)]'
$ws.Range("B2").Value = '[(-2.0, 1.0), (-2.0, 0.87), (-14.169, 0.7650108327053524), (-14.169, 0.6350108327053523), (-2.0, 0.87), (-2.0, 0.75), (-14.169, 0.6350108327053523), (-14.169, 0.5150108327053524)]'
$ws.Range("B3").Value = '[(-2.0, 2.0), (-2.0, 1.87), (-13.64, 1.7366317044100097), (-13.64, 1.6066317044100098), (-2.0, 1.87), (-2.0, 1.75), (-13.64, 1.6066317044100098), (-13.64, 1.4866317044100097)]'
$ws.Range("B4").Value = '[(-2.0, 3.0), (-2.0, 2.87), (-13.653, 2.7765873638344267), (-13.653, 2.646587363834427), (-2.0, 2.87), (-2.0, 2.75), (-13.653, 2.646587363834427), (-13.653, 2.5265873638344267)]'
$ws.Range("B5").Value = '[(-2.0, 4.0), (-2.0, 3.87), (-13.605, 3.743245028480061), (-13.605, 3.613245028480061), (-2.0, 3.87), (-2.0, 3.75), (-13.605, 3.613245028480061), (-13.605, 3.493245028480061)]'
$ws.Range("B6").Value = '[(-2.0, 5.0), (-2.0, 4.87), (-14.207, 4.733979528301885), (-14.207, 4.603979528301885), (-2.0, 4.87), (-2.0, 4.75), (-14.207, 4.603979528301885), (-14.207, 4.483979528301885)]'
$ws.Range("B7").Value = '[(-2.0, 6.0), (-2.0, 5.87), (-10.541, 5.844206061042327), (-10.541, 5.7142060610423275), (-2.0, 5.87), (-2.0, 5.75), (-10.541, 5.7142060610423275), (-10.541, 5.594206061042327)]'
$ws.Range("B8").Value = '[(-2.0, 7.0), (-2.0, 6.87), (-10.654, 6.890557444984277), (-10.654, 6.760557444984277), (-2.0, 6.87), (-2.0, 6.75), (-10.654, 6.760557444984277), (-10.654, 6.640557444984277)]'
$ws.Range("B9").Value = '[(-2.0, 8.0), (-2.0, 7.87), (-10.673, 7.862074598814507), (-10.673, 7.732074598814507), (-2.0, 7.87), (-2.0, 7.75), (-10.673, 7.732074598814507), (-10.673, 7.612074598814507)]'
$ws.Range("B10").Value = '[(-2.0, 9.0), (-2.0, 8.87), (-10.704, 8.906585699841107), (-10.704, 8.776585699841107), (-2.0, 8.87), (-2.0, 8.75), (-10.704, 8.776585699841107), (-10.704, 8.656585699841107)]'
$ws.Range("B11").Value = '[(-2.0, 10.0), (-2.0, 9.87), (-10.767, 9.874775122021248), (-10.767, 9.744775122021247), (-2.0, 9.87), (-2.0, 9.75), (-10.767, 9.744775122021247), (-10.767, 9.624775122021248)]'
$ws.Range("B12").Value = '[(-2.0, 11.0), (-2.0, 10.87), (-10.728, 10.87606615936827), (-10.728, 10.746066159368269), (-2.0, 10.87), (-2.0, 10.75), (-10.728, 10.746066159368269), (-10.728, 10.62606615936827)]'
$ws.Range("B13").Value = '[(-2.0, 12.0), (-2.0, 11.87), (-10.747, 11.869367069138994), (-10.747, 11.739367069138993), (-2.0, 11.87), (-2.0, 11.75), (-10.747, 11.739367069138993), (-10.747, 11.619367069138994)]'
$ws.Range("B14").Value = '[(-2.0, 13.0), (-2.0, 12.87), (-10.712, 12.859261918437683), (-10.712, 12.729261918437683), (-2.0, 12.87), (-2.0, 12.75), (-10.712, 12.729261918437683), (-10.712, 12.609261918437683)]'
$ws.Range("B15").Value = '[(-2.0, 14.0), (-2.0, 13.87), (-10.721, 13.842581227436822), (-10.721, 13.712581227436822), (-2.0, 13.87), (-2.0, 13.75), (-10.721, 13.712581227436822), (-10.721, 13.592581227436822)]'
$ws.Range("B16").Value = '[(-2.0, 15.0), (-2.0, 14.87), (-10.711, 14.783197349848766), (-10.711, 14.653197349848766), (-2.0, 14.87), (-2.0, 14.75), (-10.711, 14.653197349848766), (-10.711, 14.533197349848766)]'
$ws.Range("B17").Value = '[(-2.00011, 16.0), (-2.00011, 15.87), (-10.664, 15.83059988423468), (-10.664, 15.700599884234679), (-2.00011, 15.87), (-2.00011, 15.75), (-10.664, 15.700599884234679), (-10.664, 15.58059988423468)]'
$ws.Range("B18").Value = '[(-1.99961, 17.0), (-1.99961, 16.87), (-10.548, 16.720150002575476), (-10.548, 16.590150002575477), (-1.99961, 16.87), (-1.99961, 16.75), (-10.548, 16.590150002575477), (-10.548, 16.470150002575476)]'
$ws.Range("B19").Value = '[(-1.99935, 18.0), (-1.99935, 17.87), (-10.299, 17.69027787646722), (-10.299, 17.56027787646722), (-1.99935, 17.87), (-1.99935, 17.75), (-10.299, 17.56027787646722), (-10.299, 17.44027787646722)]'
$ws.Range("B20").Value = '[(-1.99929, 19.0), (-1.99929, 18.87), (-10.014, 18.636504566666662), (-10.014, 18.506504566666663), (-1.99929, 18.87), (-1.99929, 18.75), (-10.014, 18.506504566666663), (-10.014, 18.386504566666662)]'
$ws.Range("B21").Value = '[(-1.99955, 20.0), (-1.99955, 19.87), (-9.792, 19.616670637574384), (-9.792, 19.486670637574385), (-1.99955, 19.87), (-1.99955, 19.75), (-9.792, 19.486670637574385), (-9.792, 19.366670637574384)]'
$ws.Range("B22").Value = '[(-2.0, 21.0), (-2.0, 20.87), (-9.625, 20.581407364284676), (-9.625, 20.451407364284677), (-2.0, 20.87), (-2.0, 20.75), (-9.625, 20.451407364284677), (-9.625, 20.331407364284676)]'
$ws.Range("B23").Value = '[(-1.717, 22.0), (-1.717, 21.87), (-9.685, 21.5227168), (-9.685, 21.392716800000002), (-1.717, 21.87), (-1.717, 21.75), (-9.685, 21.392716800000002), (-9.685, 21.2727168)]'
$ws.Range("B24").Value = '[(-1.4335, 23.0), (-1.4335, 22.87), (-12.241, 22.2975125), (-12.241, 22.1675125), (-1.4335, 22.87), (-1.4335, 22.75), (-12.241, 22.1675125), (-12.241, 22.0475125)]'
$ws.Range("B25").Value = '[(-1.15, 24.0), (-1.15, 23.87), (-10.958, 23.31344), (-10.958, 23.18344), (-1.15, 23.87), (-1.15, 23.75), (-10.958, 23.18344), (-10.958, 23.06344)]'
$ws.Range("B26").Value = '[(-0.8665, 25.0), (-0.8665, 24.87), (-9.651, 24.385085), (-9.651, 24.255085), (-0.8665, 24.87), (-0.8665, 24.75), (-9.651, 24.255085), (-9.651, 24.135085)]'
$ws.Range("B27").Value = '[(-0.5835, 26.0), (-0.5835, 25.87), (-8.446, 25.449625), (-8.446, 25.319625000000002), (-0.5835, 25.87), (-0.5835, 25.75), (-8.446, 25.319625000000002), (-8.446, 25.199625)]'
$ws.Range("B28").Value = '[(-0.3, 27.0), (-0.3, 26.87), (-6.71, 26.5513), (-6.71, 26.421300000000002), (-0.3, 26.87), (-0.3, 26.75), (-6.71, 26.421300000000002), (-6.71, 26.3013)]'
$ws.Range("B29").Value = '[(-0.3, 28.0), (-0.3, 27.87), (-6.105, 27.59365), (-6.105, 27.46365), (-0.3, 27.87), (-0.3, 27.75), (-6.105, 27.46365), (-6.105, 27.34365)]'
$ws.Range("B30").Value = '[(-0.3, 29.0), (-0.3, 28.87), (-5.902, 28.60786), (-5.902, 28.47786), (-0.3, 28.87), (-0.3, 28.75), (-5.902, 28.47786), (-5.902, 28.35786)]'
$ws.Range("B96").Value = '[(-0.3, 95.0), (-0.3, 94.87), (-2.48, 94.8692), (-2.48, 94.73920000000001), (-0.3, 94.87), (-0.3, 94.75), (-2.48, 94.73920000000001), (-2.48, 94.6192)]'
$ws.Range("B156").Value = '[(-0.3, 155.0), (-0.3, 154.86), (-2.846, 154.93635), (-2.846, 154.79635000000002), (-2.846, 154.93635), (-2.846, 154.80635), (-10.83, 154.73675), (-10.83, 154.60675), (-0.3, 154.86), (-0.3, 154.73), (-2.846, 154.79635000000002), (-2.846, 154.68635), (-2.846, 154.80635), (-2.846, 154.68635), (-10.83, 154.60675), (-10.83, 154.48675)]'
$ws.Range("B159").Value = '[(-0.3, 158.0), (-0.3, 157.87), (-3.906, 157.90985), (-3.906, 157.77985), (-3.906, 157.90985), (-3.906, 157.77985), (-13.443, 157.671425), (-13.443, 157.541425), (-0.3, 157.87), (-0.3, 157.75), (-3.906, 157.77985), (-3.906, 157.65985), (-3.906, 157.77985), (-3.906, 157.65985), (-13.443, 157.541425), (-13.443, 157.421425)]'
$ws.Range("B160").Value = '[(-0.3, 159.0), (-0.3, 158.89), (-4.157, 158.903575), (-4.157, 158.79357499999998), (-4.157, 158.903575), (-4.157, 158.773575), (-14.083, 158.655425), (-14.083, 158.525425), (-0.3, 158.89), (-0.3, 158.68), (-4.157, 158.79357499999998), (-4.157, 158.653575), (-4.157, 158.773575), (-4.157, 158.653575), (-14.083, 158.525425), (-14.083, 158.405425)]'
$ws.Range("B163").Value = '[(-0.3, 162.0), (-0.3, 161.89), (-4.481, 161.895475), (-4.481, 161.785475), (-4.481, 161.895475), (-4.481, 161.765475), (-11.438, 161.72155), (-11.438, 161.59155), (-0.3, 161.89), (-0.3, 161.62), (-4.481, 161.785475), (-4.481, 161.645475), (-4.481, 161.765475), (-4.481, 161.645475), (-11.438, 161.59155), (-11.438, 161.47155)]'
$ws.Range("B164").Value = '[(-0.3, 163.0), (-0.3, 162.89), (-4.635, 162.891625), (-4.635, 162.781625), (-4.635, 162.891625), (-4.635, 162.761625), (-11.526, 162.71935), (-11.526, 162.58935), (-0.3, 162.89), (-0.3, 162.62), (-4.635, 162.781625), (-4.635, 162.641625), (-4.635, 162.761625), (-4.635, 162.641625), (-11.526, 162.58935), (-11.526, 162.46935)]'
$ws.Range("B165").Value = '[(-0.3, 164.0), (-0.3, 163.88), (-4.568, 163.8933), (-4.568, 163.7733), (-4.568, 163.8933), (-4.568, 163.76330000000002), (-11.542, 163.71895), (-11.542, 163.58895), (-0.3, 163.88), (-0.3, 163.58), (-4.568, 163.7733), (-4.568, 163.6433), (-4.568, 163.76330000000002), (-4.568, 163.6433), (-11.542, 163.58895), (-11.542, 163.46895)]'
$ws.Range("B166").Value = '[(-0.3, 165.0), (-0.3, 164.88), (-4.461, 164.895975), (-4.461, 164.775975), (-4.461, 164.895975), (-4.461, 164.765975), (-11.467, 164.720825), (-11.467, 164.590825), (-0.3, 164.88), (-0.3, 164.58), (-4.461, 164.775975), (-4.461, 164.645975), (-4.461, 164.765975), (-4.461, 164.645975), (-11.467, 164.590825), (-11.467, 164.470825)]'
$ws.Range("B167").Value = '[(-0.3, 166.0), (-0.3, 165.89), (-4.253, 165.901175), (-4.253, 165.79117499999998), (-4.253, 165.901175), (-4.253, 165.771175), (-11.362, 165.72345), (-11.362, 165.59345000000002), (-0.3, 165.89), (-0.3, 165.68), (-4.253, 165.79117499999998), (-4.253, 165.651175), (-4.253, 165.771175), (-4.253, 165.651175), (-11.362, 165.59345000000002), (-11.362, 165.47345)]'
$ws.Range("B168").Value = '[(-0.3, 167.0), (-0.3, 166.87), (-4.252, 166.9012), (-4.252, 166.7712), (-4.252, 166.9012), (-4.252, 166.7712), (-11.259, 166.726025), (-11.259, 166.596025), (-0.3, 166.87), (-0.3, 166.75), (-4.252, 166.7712), (-4.252, 166.6512), (-4.252, 166.7712), (-4.252, 166.6512), (-11.259, 166.596025), (-11.259, 166.476025)]'
$ws.Range("B169").Value = '[(-0.3, 168.0), (-0.3, 167.86), (-4.164, 167.9034), (-4.164, 167.76340000000002), (-4.164, 167.9034), (-4.164, 167.7734), (-11.148, 167.7288), (-11.148, 167.5988), (-0.3, 167.86), (-0.3, 167.74), (-4.164, 167.76340000000002), (-4.164, 167.6534), (-4.164, 167.7734), (-4.164, 167.6534), (-11.148, 167.5988), (-11.148, 167.4788)]'
$ws.Range("B170").Value = '[(-0.3, 169.0), (-0.3, 168.87), (-4.036, 168.9066), (-4.036, 168.7766), (-4.036, 168.9066), (-4.036, 168.7766), (-13.582, 168.66795), (-13.582, 168.53795), (-0.3, 168.87), (-0.3, 168.75), (-4.036, 168.7766), (-4.036, 168.6566), (-4.036, 168.7766), (-4.036, 168.6566), (-13.582, 168.53795), (-13.582, 168.41795)]'
$ws.Range("B171").Value = '[(-0.3, 170.0), (-0.3, 169.87), (-4.007, 169.907325), (-4.007, 169.777325), (-4.007, 169.907325), (-4.007, 169.777325), (-13.19, 169.67775), (-13.19, 169.54775), (-0.3, 169.87), (-0.3, 169.75), (-4.007, 169.777325), (-4.007, 169.657325), (-4.007, 169.777325), (-4.007, 169.657325), (-13.19, 169.54775), (-13.19, 169.42775)]'
$ws.Range("B172").Value = '[(-0.3, 171.0), (-0.3, 170.89), (-3.944, 170.9089), (-3.944, 170.79889999999997), (-3.944, 170.9089), (-3.944, 170.7789), (-12.833, 170.686675), (-12.833, 170.556675), (-0.3, 170.89), (-0.3, 170.67), (-3.944, 170.79889999999997), (-3.944, 170.6589), (-3.944, 170.7789), (-3.944, 170.6589), (-12.833, 170.556675), (-12.833, 170.436675)]'
$ws.Range("B173").Value = '[(-0.3, 172.0), (-0.3, 171.89), (-3.841, 171.911475), (-3.841, 171.80147499999998), (-3.841, 171.911475), (-3.841, 171.781475), (-12.411, 171.697225), (-12.411, 171.567225), (-0.3, 171.89), (-0.3, 171.66), (-3.841, 171.80147499999998), (-3.841, 171.661475), (-3.841, 171.781475), (-3.841, 171.661475), (-12.411, 171.567225), (-12.411, 171.447225)]'
$ws.Range("B175").Value = '[(-0.3, 174.0), (-0.3, 173.86), (-4.002, 173.88920942046707), (-4.002, 173.74920942046708), (-4.002, 173.88920942046707), (-4.002, 173.75920942046707), (-11.644, 173.66050558232808), (-11.644, 173.5305055823281), (-0.3, 173.86), (-0.3, 173.73), (-4.002, 173.74920942046708), (-4.002, 173.63920942046707), (-4.002, 173.75920942046707), (-4.002, 173.63920942046707), (-11.644, 173.5305055823281), (-11.644, 173.41050558232808)]'
$ws.Range("B176").Value = '[(-0.3, 175.0), (-0.3, 174.87), (-4.157, 174.86302894558347), (-4.157, 174.73302894558347), (-4.157, 174.86302894558347), (-4.157, 174.73302894558347), (-11.37, 174.60687851376954), (-11.37, 174.47687851376955), (-0.3, 174.87), (-0.3, 174.75), (-4.157, 174.73302894558347), (-4.157, 174.61302894558347), (-4.157, 174.73302894558347), (-4.157, 174.61302894558347), (-11.37, 174.47687851376955), (-11.37, 174.35687851376954)]'
$ws.Range("B180").Value = '[(-0.3, 179.0), (-0.3, 178.89), (-3.955, 178.78854818079898), (-3.955, 178.67854818079897), (-3.955, 178.78854818079898), (-3.955, 178.658548180799), (-10.974, 178.38247969407612), (-10.974, 178.25247969407613), (-0.3, 178.89), (-0.3, 178.67), (-3.955, 178.67854818079897), (-3.955, 178.53854818079898), (-3.955, 178.658548180799), (-3.955, 178.53854818079898), (-10.974, 178.25247969407613), (-10.974, 178.13247969407612)]'
$ws.Range("B181").Value = '[(-0.3, 180.0), (-0.3, 179.89), (-3.483, 179.80902), (-3.483, 179.69902), (-3.483, 179.80902), (-3.483, 179.67902), (-10.493, 179.38842), (-10.493, 179.25842), (-0.3, 179.89), (-0.3, 179.69), (-3.483, 179.69902), (-3.483, 179.55902), (-3.483, 179.67902), (-3.483, 179.55902), (-10.493, 179.25842), (-10.493, 179.13842)]'
$ws.Range("B182").Value = '[(-0.3, 181.0), (-0.3, 180.89), (-2.824, 180.84856), (-2.824, 180.73855999999998), (-2.824, 180.84856), (-2.824, 180.71856), (-9.811, 180.42934), (-9.811, 180.29934), (-0.3, 180.89), (-0.3, 180.66), (-2.824, 180.73855999999998), (-2.824, 180.59856), (-2.824, 180.71856), (-2.824, 180.59856), (-9.811, 180.29934), (-9.811, 180.17934)]'
$ws.Range("B183").Value = '[(-0.3, 182.0), (-0.3, 181.89), (-9.17, 181.4678), (-9.17, 181.3578), (-0.3, 181.89), (-0.3, 181.67), (-9.17, 181.3578), (-9.17, 181.2178)]'
$ws.Range("B184").Value = '[(-0.3, 183.0), (-0.3, 182.89), (-8.644, 182.49936), (-8.644, 182.38935999999998), (-0.3, 182.89), (-0.3, 182.63), (-8.644, 182.38935999999998), (-8.644, 182.24936)]'
$ws.Range("B185").Value = '[(-0.3, 184.0), (-0.3, 183.89), (-7.857, 183.54658), (-7.857, 183.43658), (-0.3, 183.89), (-0.3, 183.61), (-7.857, 183.43658), (-7.857, 183.29658)]'
$ws.Range("B186").Value = '[(-0.3, 185.0), (-0.3, 184.89), (-6.888, 184.60472), (-6.888, 184.49471999999997), (-0.3, 184.89), (-0.3, 184.63), (-6.888, 184.49471999999997), (-6.888, 184.35472)]'
$ws.Range("B187").Value = '[(-0.3, 186.0), (-0.3, 185.89), (-5.892, 185.66448), (-5.892, 185.55447999999998), (-0.3, 185.89), (-0.3, 185.68), (-5.892, 185.55447999999998), (-5.892, 185.41448)]'
$ws.Range("B188").Value = '[(-0.3, 187.0), (-0.3, 186.89), (-4.905, 186.7237), (-4.905, 186.6137), (-0.3, 186.89), (-0.3, 186.68), (-4.905, 186.6137), (-4.905, 186.4737)]'
$ws.Range("B189").Value = '[(-0.3, 188.0), (-0.3, 187.89), (-4.146, 187.76924), (-4.146, 187.65923999999998), (-0.3, 187.89), (-0.3, 187.67), (-4.146, 187.65923999999998), (-4.146, 187.51924)]'
$ws.Range("B214").Value = '[(-0.3, 213.0), (-0.3, 212.89), (-2.761, 213.07269782605752), (-2.761, 212.9626978260575), (-0.3, 212.89), (-0.3, 212.75), (-2.761, 212.9626978260575), (-2.761, 212.75269782605753)]'
$ws.Range("B215").Value = '[(-0.3, 214.0), (-0.3, 213.87), (-3.359, 214.09941811729925), (-3.359, 213.96941811729926), (-0.3, 213.87), (-0.3, 213.75), (-3.359, 213.96941811729926), (-3.359, 213.84941811729925)]'
$ws.Range("B216").Value = '[(-0.3, 215.0), (-0.3, 214.87), (-4.285, 215.14130989283044), (-4.285, 215.01130989283044), (-0.3, 214.87), (-0.3, 214.75), (-4.285, 215.01130989283044), (-4.285, 214.89130989283044)]'
$ws.Range("B217").Value = '[(-0.3, 216.0), (-0.3, 215.87), (-2.315, 216.0764202606724), (-2.315, 215.9464202606724), (-2.315, 216.0764202606724), (-2.315, 215.9464202606724), (-5.644, 216.20267487495445), (-5.644, 216.07267487495446), (-0.3, 215.87), (-0.3, 215.75), (-2.315, 215.9464202606724), (-2.315, 215.8264202606724), (-2.315, 215.9464202606724), (-2.315, 215.8264202606724), (-5.644, 216.07267487495446), (-5.644, 215.95267487495445)]'
$ws.Range("B218").Value = '[(-0.3, 217.0), (-0.3, 216.87), (-4.206, 217.23903990672179), (-4.206, 217.1090399067218), (-4.206, 217.23903990672179), (-4.206, 217.12903990672177), (-7.515, 217.4415445281612), (-7.515, 217.33154452816117), (-0.3, 216.87), (-0.3, 216.75), (-4.206, 217.1090399067218), (-4.206, 216.98903990672179), (-4.206, 217.12903990672177), (-4.206, 216.98903990672179), (-7.515, 217.33154452816117), (-7.515, 217.11154452816118)]'
$ws.Range("B223").Value = '[(-0.3, 222.0), (-0.3, 221.89), (-10.766, 222.73262), (-10.766, 222.62261999999998), (-10.766, 222.73262), (-10.766, 222.62261999999998), (-14.023, 222.96061), (-14.023, 222.85061), (-0.3, 221.89), (-0.3, 221.75), (-10.766, 222.62261999999998), (-10.766, 222.42262), (-10.766, 222.62261999999998), (-10.766, 222.42262), (-14.023, 222.85061), (-14.023, 222.55061)]'
$ws.Range("B224").Value = '[(-0.3, 223.0), (-0.3, 222.89), (-11.007, 223.74949), (-11.007, 223.63949), (-11.007, 223.74949), (-11.007, 223.58949), (-14.317, 223.98119), (-14.317, 223.82119), (-0.3, 222.89), (-0.3, 222.75), (-11.007, 223.63949), (-11.007, 223.41949), (-11.007, 223.58949), (-11.007, 223.41949), (-14.317, 223.82119), (-14.317, 223.52119)]'
$ws.Range("B225").Value = '[(-0.3, 224.0), (-0.3, 223.88), (-11.092, 224.75544), (-11.092, 224.51543999999998), (-11.092, 224.75544), (-11.092, 224.51543999999998), (-14.44, 224.9898), (-14.44, 224.6598), (-0.3, 223.88), (-0.3, 223.76), (-11.092, 224.51543999999998), (-11.092, 224.27544), (-11.092, 224.51543999999998), (-11.092, 224.27544), (-14.44, 224.6598), (-14.44, 224.3298)]'
$ws.Range("B226").Value = '[(-0.3, 225.0), (-0.3, 224.88), (-10.831, 225.73717), (-10.831, 225.47717), (-10.831, 225.73717), (-10.831, 225.47717), (-14.2, 225.973), (-14.2, 225.633), (-0.3, 224.88), (-0.3, 224.76), (-10.831, 225.47717), (-10.831, 225.21716999999998), (-10.831, 225.47717), (-10.831, 225.21716999999998), (-14.2, 225.633), (-14.2, 225.293)]'
$ws.Range("B227").Value = '[(-0.3, 226.0), (-0.3, 225.88), (-10.166, 226.69062), (-10.166, 226.43062), (-10.166, 226.69062), (-10.166, 226.43062), (-13.518, 226.92526), (-13.518, 226.59526), (-0.3, 225.88), (-0.3, 225.76), (-10.166, 226.43062), (-10.166, 226.17061999999999), (-10.166, 226.43062), (-10.166, 226.17061999999999), (-13.518, 226.59526), (-13.518, 226.26526)]'
$ws.Range("B228").Value = '[(-0.3, 227.0), (-0.3, 226.88), (-9.205, 227.62335), (-9.205, 227.38334999999998), (-9.205, 227.62335), (-9.205, 227.38334999999998), (-12.504, 227.85428), (-12.504, 227.51427999999999), (-0.3, 226.88), (-0.3, 226.76), (-9.205, 227.38334999999998), (-9.205, 227.14335), (-9.205, 227.38334999999998), (-9.205, 227.14335), (-12.504, 227.51427999999999), (-12.504, 227.17427999999998)]'
$ws.Range("B229").Value = '[(-0.3, 228.0), (-0.3, 227.89), (-8.19, 228.4920039533576), (-8.19, 228.38200395335758), (-8.19, 228.4920039533576), (-8.19, 228.3120039533576), (-11.646, 228.70751290935303), (-11.646, 228.52751290935302), (-0.3, 227.89), (-0.3, 227.75), (-8.19, 228.38200395335758), (-8.19, 228.1820039533576), (-8.19, 228.3120039533576), (-8.19, 228.1820039533576), (-11.646, 228.52751290935302), (-11.646, 228.22751290935304)]'
$ws.Range("B230").Value = '[(-0.3, 229.0), (-0.3, 228.87), (-7.184, 229.37280157785662), (-7.184, 229.24280157785662), (-7.184, 229.37280157785662), (-7.184, 229.23280157785663), (-10.707, 229.56358890481607), (-10.707, 229.42358890481609), (-0.3, 228.87), (-0.3, 228.75), (-7.184, 229.24280157785662), (-7.184, 229.12280157785662), (-7.184, 229.23280157785663), (-7.184, 229.12280157785662), (-10.707, 229.42358890481609), (-10.707, 229.12358890481607)]'
$ws.Range("B233").Value = '[(-0.3, 232.0), (-0.3, 231.87), (-6.473, 232.1823838499305), (-6.473, 232.0523838499305), (-6.473, 232.1823838499305), (-6.473, 232.0423838499305), (-9.977, 232.2859109858703), (-9.977, 232.1459109858703), (-0.3, 231.87), (-0.3, 231.75), (-6.473, 232.0523838499305), (-6.473, 231.9323838499305), (-6.473, 232.0423838499305), (-6.473, 231.9323838499305), (-9.977, 232.1459109858703), (-9.977, 231.8459109858703)]'
$ws.Range("B234").Value = '[(-0.3, 233.0), (-0.3, 232.87), (-6.598, 233.07055554636707), (-6.598, 232.94055554636708), (-6.598, 233.07055554636707), (-6.598, 232.96055554636706), (-10.064, 233.10938462285299), (-10.064, 232.99938462285297), (-0.3, 232.87), (-0.3, 232.75), (-6.598, 232.94055554636708), (-6.598, 232.82055554636707), (-6.598, 232.96055554636706), (-6.598, 232.82055554636707), (-10.064, 232.99938462285297), (-10.064, 232.759384622853)]'
$ws.Range("B235").Value = '[(-0.3, 234.0), (-0.3, 233.87), (-6.694, 233.8737827975527), (-6.694, 233.7437827975527), (-6.694, 233.8737827975527), (-6.694, 233.7437827975527), (-10.232, 233.80394287539778), (-10.232, 233.67394287539778), (-0.3, 233.87), (-0.3, 233.75), (-6.694, 233.7437827975527), (-6.694, 233.6237827975527), (-6.694, 233.7437827975527), (-6.694, 233.6237827975527), (-10.232, 233.67394287539778), (-10.232, 233.55394287539778)]'
$ws.Range("B236").Value = '[(-0.3, 235.0), (-0.3, 234.87), (-6.865, 234.72144607427487), (-6.865, 234.59144607427487), (-6.865, 234.72144607427487), (-6.865, 234.59144607427487), (-10.356, 234.57332242542392), (-10.356, 234.44332242542393), (-0.3, 234.87), (-0.3, 234.75), (-6.865, 234.59144607427487), (-6.865, 234.47144607427487), (-6.865, 234.59144607427487), (-6.865, 234.47144607427487), (-10.356, 234.44332242542393), (-10.356, 234.32332242542392)]'
$ws.Range("B237").Value = '[(-0.3, 236.0), (-0.3, 235.87), (-6.433, 235.61098289848098), (-6.433, 235.480982898481), (-6.433, 235.61098289848098), (-6.433, 235.480982898481), (-10.083, 235.3794628559986), (-10.083, 235.24946285599862), (-0.3, 235.87), (-0.3, 235.75), (-6.433, 235.480982898481), (-6.433, 235.36098289848098), (-6.433, 235.480982898481), (-6.433, 235.36098289848098), (-10.083, 235.24946285599862), (-10.083, 235.1294628559986)]'
$ws.Range("B238").Value = '[(-0.3, 237.0), (-0.3, 236.87), (-5.622, 236.62746), (-5.622, 236.49746000000002), (-5.622, 236.62746), (-5.622, 236.49746000000002), (-9.234, 236.37462), (-9.234, 236.24462), (-0.3, 236.87), (-0.3, 236.75), (-5.622, 236.49746000000002), (-5.622, 236.37746), (-5.622, 236.49746000000002), (-5.622, 236.37746), (-9.234, 236.24462), (-9.234, 236.12462)]'
$ws.Range("B239").Value = '[(-0.3, 238.0), (-0.3, 237.87), (-4.904, 237.67772), (-4.904, 237.54772), (-4.904, 237.67772), (-4.904, 237.54772), (-8.564, 237.42152), (-8.564, 237.29152), (-0.3, 237.87), (-0.3, 237.75), (-4.904, 237.54772), (-4.904, 237.42772), (-4.904, 237.54772), (-4.904, 237.42772), (-8.564, 237.29152), (-8.564, 237.17152)]'
$ws.Range("B240").Value = '[(-0.3, 239.0), (-0.3, 238.87), (-4.102, 238.73386), (-4.102, 238.60386), (-4.102, 238.73386), (-4.102, 238.60386), (-7.913, 238.46709), (-7.913, 238.33709000000002), (-0.3, 238.87), (-0.3, 238.75), (-4.102, 238.60386), (-4.102, 238.48386), (-4.102, 238.60386), (-4.102, 238.48386), (-7.913, 238.33709000000002), (-7.913, 238.21709)]'
$ws.Range("B241").Value = '[(-0.3, 240.0), (-0.3, 239.87), (-3.709, 239.76137), (-3.709, 239.63137), (-3.709, 239.76137), (-3.709, 239.63137), (-7.59, 239.4897), (-7.59, 239.3597), (-0.3, 239.87), (-0.3, 239.75), (-3.709, 239.63137), (-3.709, 239.51137), (-3.709, 239.63137), (-3.709, 239.51137), (-7.59, 239.3597), (-7.59, 239.2397)]'
$ws.Range("B242").Value = '[(-0.3, 241.0), (-0.3, 240.87), (-3.488, 240.84073776727553), (-3.488, 240.71073776727553), (-3.488, 240.84073776727553), (-3.488, 240.71073776727553), (-7.284, 240.65110180886208), (-7.284, 240.52110180886208), (-0.3, 240.87), (-0.3, 240.75), (-3.488, 240.71073776727553), (-3.488, 240.59073776727553), (-3.488, 240.71073776727553), (-3.488, 240.59073776727553), (-7.284, 240.52110180886208), (-7.284, 240.40110180886208)]'
$ws.Range("B243").Value = '[(-0.3, 242.0), (-0.3, 241.87), (-3.388, 241.91058144584278), (-3.388, 241.78058144584278), (-3.388, 241.91058144584278), (-3.388, 241.78058144584278), (-7.389, 241.7947253463664), (-7.389, 241.6647253463664), (-0.3, 241.87), (-0.3, 241.75), (-3.388, 241.78058144584278), (-3.388, 241.66058144584278), (-3.388, 241.78058144584278), (-3.388, 241.66058144584278), (-7.389, 241.6647253463664), (-7.389, 241.5447253463664)]'
$ws.Range("B249").Value = '[(-0.3, 248.0), (-0.3, 247.86), (-3.123, 248.19761), (-3.123, 248.05761), (-3.123, 248.19761), (-3.123, 248.08760999999998), (-6.875, 248.46025), (-6.875, 248.35025), (-0.3, 247.86), (-0.3, 247.75), (-3.123, 248.05761), (-3.123, 247.93761), (-3.123, 248.08760999999998), (-3.123, 247.93761), (-6.875, 248.35025), (-6.875, 248.15025)]'
$ws.Range("B266").Value = '[(-0.3, 265.0), (-0.3, 264.85), (-6.557, 264.56201), (-6.557, 264.41201), (-6.557, 264.56201), (-6.557, 264.45201), (-10.352, 264.29636), (-10.352, 264.18636), (-0.3, 264.85), (-0.3, 264.55), (-6.557, 264.41201), (-6.557, 264.25201), (-6.557, 264.45201), (-6.557, 264.25201), (-10.352, 264.18636), (-10.352, 264.04636)]'
$ws.Range("B275").Value = '[(-0.3, 274.0), (-0.3, 273.86), (-5.91, 273.66646702810965), (-5.91, 273.52646702810966), (-0.3, 273.86), (-0.3, 273.74), (-5.91, 273.52646702810966), (-5.91, 273.41646702810965)]'
$ws.Range("B293").Value = '[(-0.3, 292.0), (-0.3, 291.86), (-2.517, 291.84481), (-2.517, 291.70481), (-0.3, 291.86), (-0.3, 291.73), (-2.517, 291.70481), (-2.517, 291.59481)]'
$ws.Range("B460").Value = '[(-0.3, 459.0), (-0.3, 458.89), (-8.349, 458.91663834317296), (-8.349, 458.80663834317295), (-0.3, 458.89), (-0.3, 458.68), (-8.349, 458.80663834317295), (-8.349, 458.66663834317296)]'
$ws.Range("B463").Value = '[(-0.3, 462.0), (-0.3, 461.89), (-8.916, 462.3270511115744), (-8.916, 462.2170511115744), (-0.3, 461.89), (-0.3, 461.75), (-8.916, 462.2170511115744), (-8.916, 461.9770511115744)]'
$ws.Range("B464").Value = '[(-0.3, 463.0), (-0.3, 462.89), (-8.991, 463.4009208569551), (-8.991, 463.2909208569551), (-0.3, 462.89), (-0.3, 462.75), (-8.991, 463.2909208569551), (-8.991, 463.00092085695513)]'
$ws.Range("B468").Value = '[(-0.3, 467.0), (-0.3, 466.86), (-3.104, 467.19628), (-3.104, 467.05628), (-3.104, 467.19628), (-3.104, 467.08628), (-9.819, 467.66633), (-9.819, 467.55633), (-0.3, 466.86), (-0.3, 466.75), (-3.104, 467.05628), (-3.104, 466.93628), (-3.104, 467.08628), (-3.104, 466.93628), (-9.819, 467.55633), (-9.819, 467.25633)]'
$ws.Range("B469").Value = '[(-0.3, 468.0), (-0.3, 467.87), (-3.425, 468.21875), (-3.425, 468.08875), (-3.425, 468.21875), (-3.425, 468.10875), (-10.177, 468.69139), (-10.177, 468.58139), (-0.3, 467.87), (-0.3, 467.75), (-3.425, 468.08875), (-3.425, 467.96875), (-3.425, 468.10875), (-3.425, 467.96875), (-10.177, 468.58139), (-10.177, 468.31139)]'
$ws.Range("B471").Value = '[(-0.3, 470.0), (-0.3, 469.87), (-3.89, 470.2513), (-3.89, 470.1213), (-3.89, 470.2513), (-3.89, 470.1413), (-10.655, 470.72485), (-10.655, 470.61485), (-0.3, 469.87), (-0.3, 469.75), (-3.89, 470.1213), (-3.89, 470.0013), (-3.89, 470.1413), (-3.89, 470.0013), (-10.655, 470.61485), (-10.655, 470.31485)]'
$ws.Range("B474").Value = '[(-0.3, 473.0), (-0.3, 472.86), (-4.582, 473.29974), (-4.582, 473.15974), (-4.582, 473.29974), (-4.582, 473.18974), (-11.31, 473.7707), (-11.31, 473.66069999999996), (-0.3, 472.86), (-0.3, 472.75), (-4.582, 473.15974), (-4.582, 473.02974), (-4.582, 473.18974), (-4.582, 473.02974), (-11.31, 473.66069999999996), (-11.31, 473.36069999999995)]'
$ws.Range("B478").Value = '[(-0.3, 477.0), (-0.3, 476.87), (-4.91, 477.3227), (-4.91, 477.1927), (-4.91, 477.3227), (-4.91, 477.2027), (-11.668, 477.79576), (-11.668, 477.67575999999997), (-0.3, 476.87), (-0.3, 476.75), (-4.91, 477.1927), (-4.91, 477.0727), (-4.91, 477.2027), (-4.91, 477.0727), (-11.668, 477.67575999999997), (-11.668, 477.37575999999996)]'
$ws.Range("B480").Value = '[(-0.3, 479.0), (-0.3, 478.86), (-4.907, 479.32249), (-4.907, 479.18249000000003), (-4.907, 479.32249), (-4.907, 479.20249), (-11.606, 479.79142), (-11.606, 479.67142), (-0.3, 478.86), (-0.3, 478.75), (-4.907, 479.18249000000003), (-4.907, 479.06249), (-4.907, 479.20249), (-4.907, 479.06249), (-11.606, 479.67142), (-11.606, 479.37142)]'
$ws.Range("B481").Value = '[(-0.3, 480.0), (-0.3, 479.86), (-4.838, 480.31766), (-4.838, 480.17766), (-4.838, 480.31766), (-4.838, 480.17766), (-11.598, 480.79086), (-11.598, 480.65086), (-0.3, 479.86), (-0.3, 479.75), (-4.838, 480.17766), (-4.838, 480.04766), (-4.838, 480.17766), (-4.838, 480.04766), (-11.598, 480.65086), (-11.598, 480.35086)]'
$ws.Range("B491").Value = '[(-0.3, 490.0), (-0.3, 489.69), (-5.217, 489.877075), (-5.217, 489.687075), (-5.217, 489.877075), (-5.217, 489.767075), (-12.015, 489.707125), (-12.015, 489.597125), (-0.3, 489.69), (-0.3, 489.38), (-5.217, 489.687075), (-5.217, 489.497075), (-5.217, 489.767075), (-5.217, 489.497075), (-12.015, 489.597125), (-12.015, 489.457125)]'
$ws.Range("B492").Value = '[(-0.3, 491.0), (-0.3, 490.68), (-5.197, 490.877575), (-5.197, 490.677575), (-5.197, 490.877575), (-5.197, 490.76757499999997), (-12.024, 490.7069), (-12.024, 490.5969), (-0.3, 490.68), (-0.3, 490.36), (-5.197, 490.677575), (-5.197, 490.477575), (-5.197, 490.76757499999997), (-5.197, 490.487575), (-12.024, 490.5969), (-12.024, 490.43690000000004)]'
$ws.Range("B493").Value = '[(-0.3, 492.0), (-0.3, 491.68), (-5.221, 491.876975), (-5.221, 491.676975), (-5.221, 491.876975), (-5.221, 491.766975), (-12.004, 491.7074), (-12.004, 491.5974), (-0.3, 491.68), (-0.3, 491.36), (-5.221, 491.676975), (-5.221, 491.47697500000004), (-5.221, 491.766975), (-5.221, 491.47697500000004), (-12.004, 491.5974), (-12.004, 491.4074)]'
$ws.Range("B494").Value = '[(-0.3, 493.0), (-0.3, 492.68), (-5.159, 492.878525), (-5.159, 492.67852500000004), (-5.159, 492.878525), (-5.159, 492.768525), (-11.984, 492.7079), (-11.984, 492.5979), (-0.3, 492.68), (-0.3, 492.36), (-5.159, 492.67852500000004), (-5.159, 492.47852500000005), (-5.159, 492.768525), (-5.159, 492.48852500000004), (-11.984, 492.5979), (-11.984, 492.4079)]'
$ws.Range("B495").Value = '[(-0.3, 494.0), (-0.3, 493.65), (-5.114, 493.87965), (-5.114, 493.63965), (-5.114, 493.87965), (-5.114, 493.68965000000003), (-11.935, 493.709125), (-11.935, 493.519125), (-0.3, 493.65), (-0.3, 493.3), (-5.114, 493.63965), (-5.114, 493.39965), (-5.114, 493.68965000000003), (-5.114, 493.38965), (-11.935, 493.519125), (-11.935, 493.369125)]'
$ws.Range("B496").Value = '[(-0.3, 495.0), (-0.3, 494.66), (-5.078, 494.88055), (-5.078, 494.62055000000004), (-5.078, 494.88055), (-5.078, 494.61055000000005), (-11.886, 494.71035), (-11.886, 494.44035), (-0.3, 494.66), (-0.3, 494.32), (-5.078, 494.62055000000004), (-5.078, 494.36055000000005), (-5.078, 494.61055000000005), (-5.078, 494.35055000000006), (-11.886, 494.44035), (-11.886, 494.28035)]'
$ws.Range("B497").Value = '[(-0.3, 496.0), (-0.3, 495.78), (-5.063, 495.880925), (-5.063, 495.66092499999996), (-5.063, 495.880925), (-5.063, 495.680925), (-11.947, 495.708825), (-11.947, 495.508825), (-0.3, 495.78), (-0.3, 495.56), (-5.063, 495.66092499999996), (-5.063, 495.440925), (-5.063, 495.680925), (-5.063, 495.480925), (-11.947, 495.508825), (-11.947, 495.308825), (-0.3, 495.56), (-0.3, 495.36), (-5.063, 495.440925), (-5.063, 495.300925), (-5.063, 495.480925), (-5.063, 495.300925), (-11.947, 495.308825), (-11.947, 495.16882499999997)]'
$ws.Range("B498").Value = '[(-0.3, 497.0), (-0.3, 496.79), (-5.088, 496.8522459061705), (-5.088, 496.6422459061705), (-5.088, 496.8522459061705), (-5.088, 496.6622459061705), (-11.96, 496.6401811332389), (-11.96, 496.4501811332389), (-0.3, 496.79), (-0.3, 496.58), (-5.088, 496.6422459061705), (-5.088, 496.43224590617046), (-5.088, 496.6622459061705), (-5.088, 496.4722459061705), (-11.96, 496.4501811332389), (-11.96, 496.2601811332389), (-0.3, 496.58), (-0.3, 496.38), (-5.088, 496.43224590617046), (-5.088, 496.3222459061705), (-5.088, 496.4722459061705), (-5.088, 496.3222459061705), (-11.96, 496.2601811332389), (-11.96, 496.0801811332389)]'
$ws.Range("B499").Value = '[(-0.3, 498.0), (-0.3, 497.79), (-5.187, 497.80034608329106), (-5.187, 497.5903460832911), (-5.187, 497.80034608329106), (-5.187, 497.61034608329106), (-12.125, 497.51690043685636), (-12.125, 497.32690043685636), (-0.3, 497.79), (-0.3, 497.58), (-5.187, 497.5903460832911), (-5.187, 497.38034608329104), (-5.187, 497.61034608329106), (-5.187, 497.42034608329107), (-12.125, 497.32690043685636), (-12.125, 497.13690043685637), (-0.3, 497.58), (-0.3, 497.38), (-5.187, 497.38034608329104), (-5.187, 497.23034608329107), (-5.187, 497.42034608329107), (-5.187, 497.23034608329107), (-12.125, 497.13690043685637), (-12.125, 496.9869004368564)]'
$ws.Range("B500").Value = '[(-0.3, 499.0), (-0.3, 498.64), (-5.357, 498.74285700980244), (-5.357, 498.4528570098024), (-5.357, 498.74285700980244), (-5.357, 498.4528570098024), (-12.301, 498.389762106909), (-12.301, 498.099762106909), (-0.3, 498.64), (-0.3, 498.28), (-5.357, 498.4528570098024), (-5.357, 498.16285700980245), (-5.357, 498.4528570098024), (-5.357, 498.16285700980245), (-12.301, 498.099762106909), (-12.301, 497.909762106909)]'
$ws.Range("B501").Value = '[(-0.3, 500.0), (-0.3, 499.63), (-5.533, 499.68160462999606), (-5.533, 499.4016046299961), (-5.533, 499.68160462999606), (-5.533, 499.4016046299961), (-12.507, 499.2572802825075), (-12.507, 499.0772802825075), (-0.3, 499.63), (-0.3, 499.26), (-5.533, 499.4016046299961), (-5.533, 499.12160462999606), (-5.533, 499.4016046299961), (-5.533, 499.12160462999606), (-12.507, 499.0772802825075), (-12.507, 498.8972802825075)]'
$ws.Range("B502").Value = '[(-0.3, 501.0), (-0.3, 500.62), (-5.774, 500.61682), (-5.774, 500.35682), (-5.774, 500.61682), (-5.774, 500.38682), (-12.784, 500.12612), (-12.784, 499.89612), (-0.3, 500.62), (-0.3, 500.24), (-5.774, 500.35682), (-5.774, 500.09682000000004), (-5.774, 500.38682), (-5.774, 500.08682000000005), (-12.784, 499.89612), (-12.784, 499.76612)]'
$ws.Range("B503").Value = '[(-0.3, 502.0), (-0.3, 501.66), (-5.857, 501.61101), (-5.857, 501.39101), (-5.857, 501.61101), (-5.857, 501.39101), (-12.779, 501.12647), (-12.779, 501.00647), (-0.3, 501.66), (-0.3, 501.32), (-5.857, 501.39101), (-5.857, 501.17101), (-5.857, 501.39101), (-5.857, 501.17101), (-12.779, 501.00647), (-12.779, 500.88647)]'
$ws.Range("B504").Value = '[(-0.3, 503.0), (-0.3, 502.66), (-5.904, 502.60772), (-5.904, 502.35772), (-5.904, 502.60772), (-5.904, 502.35772), (-12.805, 502.12465), (-12.805, 501.98465), (-0.3, 502.66), (-0.3, 502.32), (-5.904, 502.35772), (-5.904, 502.10772), (-5.904, 502.35772), (-5.904, 502.10772), (-12.805, 501.98465), (-12.805, 501.84465)]'
$ws.Range("B505").Value = '[(-0.3, 504.0), (-0.3, 503.6), (-5.981, 503.60233), (-5.981, 503.32233), (-5.981, 503.60233), (-5.981, 503.31233), (-12.864, 503.12052), (-12.864, 502.83052), (-0.3, 503.6), (-0.3, 503.2), (-5.981, 503.32233), (-5.981, 503.04233), (-5.981, 503.31233), (-5.981, 503.03233), (-12.864, 502.83052), (-12.864, 502.70052)]'
$ws.Range("B506").Value = '[(-0.3, 505.0), (-0.3, 504.6), (-5.857, 504.61101), (-5.857, 504.30101), (-5.857, 504.61101), (-5.857, 504.30101), (-12.727, 504.13011), (-12.727, 503.91011), (-0.3, 504.6), (-0.3, 504.2), (-5.857, 504.30101), (-5.857, 503.99101), (-5.857, 504.30101), (-5.857, 503.99101), (-12.727, 503.91011), (-12.727, 503.69011)]'
$ws.Range("B507").Value = '[(-0.3, 506.0), (-0.3, 505.62), (-5.566, 505.63138), (-5.566, 505.37138), (-5.566, 505.63138), (-5.566, 505.40137999999996), (-12.448, 505.14964), (-12.448, 504.91963999999996), (-0.3, 505.62), (-0.3, 505.24), (-5.566, 505.37138), (-5.566, 505.11138), (-5.566, 505.40137999999996), (-5.566, 505.10138), (-12.448, 504.91963999999996), (-12.448, 504.80964)]'
$ws.Range("B508").Value = '[(-0.3, 507.0), (-0.3, 506.71), (-5.221, 506.65559939136335), (-5.221, 506.36559939136333), (-5.221, 506.65559939136335), (-5.221, 506.54559939136334), (-12.118, 506.1729066464402), (-12.118, 506.0629066464402), (-0.3, 506.71), (-0.3, 506.42), (-5.221, 506.36559939136333), (-5.221, 506.2555993913634), (-5.221, 506.54559939136334), (-5.221, 506.2555993913634), (-12.118, 506.0629066464402), (-12.118, 505.9229066464402)]'
$ws.Range("B509").Value = '[(-0.3, 508.0), (-0.3, 507.74), (-4.862, 507.7263207648414), (-4.862, 507.4663207648414), (-4.862, 507.7263207648414), (-4.862, 507.61632076484136), (-11.663, 507.3183215368024), (-11.663, 507.2083215368024), (-0.3, 507.74), (-0.3, 507.48), (-4.862, 507.4663207648414), (-4.862, 507.3263207648414), (-4.862, 507.61632076484136), (-4.862, 507.3263207648414), (-11.663, 507.2083215368024), (-11.663, 507.0683215368024)]'
$ws.Range("B510").Value = '[(-0.3, 509.0), (-0.3, 508.81), (-4.666, 508.7817164616033), (-4.666, 508.5917164616033), (-4.666, 508.7817164616033), (-4.666, 508.6717164616033), (-11.551, 508.43749242086557), (-11.551, 508.32749242086555), (-0.3, 508.81), (-0.3, 508.51), (-4.666, 508.5917164616033), (-4.666, 508.4317164616033), (-4.666, 508.6717164616033), (-4.666, 508.4317164616033), (-11.551, 508.32749242086555), (-11.551, 508.18749242086557)]'
$ws.Range("B511").Value = '[(-0.3, 510.0), (-0.3, 509.79), (-4.571, 509.8291540422877), (-4.571, 509.61915404228773), (-4.571, 509.8291540422877), (-4.571, 509.7191540422877), (-11.467, 509.5533044229049), (-11.467, 509.4433044229049), (-0.3, 509.79), (-0.3, 509.58), (-4.571, 509.61915404228773), (-4.571, 509.5091540422877), (-4.571, 509.7191540422877), (-4.571, 509.5091540422877), (-11.467, 509.4433044229049), (-11.467, 509.3033044229049)]'
$ws.Range("B512").Value = '[(-0.3, 511.0), (-0.3, 510.89), (-4.512, 510.8736123682158), (-4.512, 510.7636123682158), (-4.512, 510.8736123682158), (-4.512, 510.7336123682158), (-11.402, 510.666867168075), (-11.402, 510.526867168075), (-0.3, 510.89), (-0.3, 510.65), (-4.512, 510.7636123682158), (-4.512, 510.6136123682158), (-4.512, 510.7336123682158), (-4.512, 510.6136123682158), (-11.402, 510.526867168075), (-11.402, 510.416867168075)]'
$ws.Range("B513").Value = '[(-0.3, 512.0), (-0.3, 511.8), (-4.393, 511.9270935479339), (-4.393, 511.7270935479339), (-4.393, 511.9270935479339), (-4.393, 511.7470935479339), (-11.313, 511.80383123464367), (-11.313, 511.62383123464366), (-0.3, 511.8), (-0.3, 511.6), (-4.393, 511.7270935479339), (-4.393, 511.6170935479339), (-4.393, 511.7470935479339), (-4.393, 511.6170935479339), (-11.313, 511.62383123464366), (-11.313, 511.32383123464365)]'
$ws.Range("B514").Value = '[(-0.3, 513.0), (-0.3, 512.75), (-4.237, 512.9865706123701), (-4.237, 512.7365706123701), (-4.237, 512.9865706123701), (-4.237, 512.7665706123701), (-11.169, 512.962925066256), (-11.169, 512.642925066256), (-0.3, 512.75), (-0.3, 512.5), (-4.237, 512.7365706123701), (-4.237, 512.54657061237), (-4.237, 512.7665706123701), (-4.237, 512.54657061237), (-11.169, 512.642925066256), (-11.169, 512.3229250662561)]'
$ws.Range("B515").Value = '[(-0.3, 514.0), (-0.3, 513.72), (-4.16, 514.0424226755699), (-4.16, 513.76242267557), (-4.16, 514.0424226755699), (-4.16, 513.8024226755699), (-11.078, 514.1184537816821), (-11.078, 513.788453781682), (-0.3, 513.72), (-0.3, 513.45), (-4.16, 513.76242267557), (-4.16, 513.5624226755699), (-4.16, 513.8024226755699), (-4.16, 513.5624226755699), (-11.078, 513.788453781682), (-11.078, 513.4584537816821)]'
$ws.Range("B516").Value = '[(-0.3, 515.0), (-0.3, 514.85), (-4.027, 515.0942611932987), (-4.027, 514.9442611932988), (-4.027, 515.0942611932987), (-4.027, 514.8042611932988), (-10.92, 515.2685950825953), (-10.92, 514.9785950825953), (-0.3, 514.85), (-0.3, 514.7), (-4.027, 514.9442611932988), (-4.027, 514.7942611932988), (-4.027, 514.8042611932988), (-4.027, 514.6542611932987), (-10.92, 514.9785950825953), (-10.92, 514.6985950825953), (-0.3, 514.7), (-0.3, 514.57), (-4.027, 514.7942611932988), (-4.027, 514.6542611932987), (-4.027, 514.6542611932987), (-4.027, 514.6542611932987), (-10.92, 514.6985950825953), (-10.92, 514.6985950825953)]'
$ws.Range("B517").Value = '[(-0.3, 516.0), (-0.3, 515.85), (-3.916, 516.1301967007082), (-3.916, 515.9801967007082), (-3.916, 516.1301967007082), (-3.916, 515.8601967007082), (-10.782, 516.377412006865), (-10.782, 516.107412006865), (-0.3, 515.85), (-0.3, 515.7), (-3.916, 515.9801967007082), (-3.916, 515.8301967007083), (-3.916, 515.8601967007082), (-3.916, 515.7201967007082), (-10.782, 516.107412006865), (-10.782, 515.847412006865), (-0.3, 515.7), (-0.3, 515.58), (-3.916, 515.8301967007083), (-3.916, 515.7201967007082), (-3.916, 515.7201967007082), (-3.916, 515.7201967007082), (-10.782, 515.847412006865), (-10.782, 515.847412006865)]'
$ws.Range("B518").Value = '[(-0.3, 517.0), (-0.3, 516.79), (-4.131, 517.1789843601017), (-4.131, 516.9689843601017), (-4.131, 517.1789843601017), (-4.131, 516.9189843601017), (-10.953, 517.4977082715121), (-10.953, 517.2377082715121), (-0.3, 516.79), (-0.3, 516.58), (-4.131, 516.9689843601017), (-4.131, 516.7989843601017), (-4.131, 516.9189843601017), (-4.131, 516.7989843601017), (-10.953, 517.2377082715121), (-10.953, 516.9877082715121)]'
$ws.Range("B519").Value = '[(-0.3, 518.0), (-0.3, 517.76), (-4.733, 518.2546062349746), (-4.733, 518.0146062349746), (-4.733, 518.2546062349746), (-4.733, 518.0146062349746), (-11.574, 518.6475142551554), (-11.574, 518.3275142551554), (-0.3, 517.76), (-0.3, 517.58), (-4.733, 518.0146062349746), (-4.733, 517.7846062349746), (-4.733, 518.0146062349746), (-4.733, 517.7746062349746), (-11.574, 518.3275142551554), (-11.574, 518.0075142551555)]'
$ws.Range("B520").Value = '[(-0.3, 519.0), (-0.3, 518.74), (-5.234, 519.3362451030761), (-5.234, 519.0762451030761), (-5.234, 519.3362451030761), (-5.234, 519.0762451030761), (-12.18, 519.8096051529274), (-12.18, 519.4596051529273), (-0.3, 518.74), (-0.3, 518.59), (-5.234, 519.0762451030761), (-5.234, 518.8262451030761), (-5.234, 519.0762451030761), (-5.234, 518.8162451030761), (-12.18, 519.4596051529273), (-12.18, 519.1096051529273)]'
$ws.Range("B521").Value = '[(-0.3, 520.0), (-0.3, 519.73), (-5.555, 520.36785), (-5.555, 520.09785), (-5.555, 520.36785), (-5.555, 520.09785), (-12.458, 520.85106), (-12.458, 520.5110599999999), (-0.3, 519.73), (-0.3, 519.55), (-5.555, 520.09785), (-5.555, 519.82785), (-5.555, 520.09785), (-5.555, 519.82785), (-12.458, 520.5110599999999), (-12.458, 520.17106)]'
$ws.Range("B522").Value = '[(-0.3, 521.0), (-0.3, 520.75), (-5.537, 521.36659), (-5.537, 521.11659), (-5.537, 521.36659), (-5.537, 521.11659), (-12.466, 521.85162), (-12.466, 521.51162), (-0.3, 520.75), (-0.3, 520.6), (-5.537, 521.11659), (-5.537, 520.86659), (-5.537, 521.11659), (-5.537, 520.86659), (-12.466, 521.51162), (-12.466, 521.1716200000001)]'
$ws.Range("B533").Value = '[(-0.3, 532.0), (-0.3, 531.89), (-7.239, 531.8360117809314), (-7.239, 531.7260117809313), (-0.3, 531.89), (-0.3, 531.66), (-7.239, 531.7260117809313), (-7.239, 531.5860117809314)]'
$ws.Range("B537").Value = '[(-0.3, 536.0), (-0.3, 535.87), (-7.123, 535.829425), (-7.123, 535.699425), (-0.3, 535.87), (-0.3, 535.75), (-7.123, 535.699425), (-7.123, 535.579425)]'
$ws.Range("B542").Value = '[(-0.3, 541.0), (-0.3, 540.87), (-6.828, 540.8368), (-6.828, 540.7068), (-0.3, 540.87), (-0.3, 540.75), (-6.828, 540.7068), (-6.828, 540.5868)]'
$ws.Range("B543").Value = '[(-0.3, 542.0), (-0.3, 541.87), (-6.768, 541.8383), (-6.768, 541.7083), (-0.3, 541.87), (-0.3, 541.75), (-6.768, 541.7083), (-6.768, 541.5883)]'
$ws.Range("B544").Value = '[(-0.3, 543.0), (-0.3, 542.87), (-6.686, 542.84035), (-6.686, 542.71035), (-0.3, 542.87), (-0.3, 542.75), (-6.686, 542.71035), (-6.686, 542.59035)]'
$ws.Range("B545").Value = '[(-0.3, 544.0), (-0.3, 543.87), (-6.632, 543.8417), (-6.632, 543.7117), (-0.3, 543.87), (-0.3, 543.75), (-6.632, 543.7117), (-6.632, 543.5917)]'
$ws.Range("B546").Value = '[(-0.3, 545.0), (-0.3, 544.87), (-6.632, 544.8417), (-6.632, 544.7117), (-0.3, 544.87), (-0.3, 544.75), (-6.632, 544.7117), (-6.632, 544.5917)]'
$ws.Range("B547").Value = '[(-0.3, 546.0), (-0.3, 545.87), (-6.624, 545.8419), (-6.624, 545.7119), (-0.3, 545.87), (-0.3, 545.75), (-6.624, 545.7119), (-6.624, 545.5919)]'
$ws.Range("B548").Value = '[(-0.3, 547.0), (-0.3, 546.87), (-6.615, 546.842125), (-6.615, 546.712125), (-0.3, 546.87), (-0.3, 546.75), (-6.615, 546.712125), (-6.615, 546.592125)]'
$ws.Range("B549").Value = '[(-0.3, 548.0), (-0.3, 547.87), (-6.603, 547.842425), (-6.603, 547.712425), (-0.3, 547.87), (-0.3, 547.75), (-6.603, 547.712425), (-6.603, 547.592425)]'
$ws.Range("B550").Value = '[(-0.3, 549.0), (-0.3, 548.87), (-6.588, 548.8428), (-6.588, 548.7128), (-0.3, 548.87), (-0.3, 548.75), (-6.588, 548.7128), (-6.588, 548.5928)]'
$ws.Range("B551").Value = '[(-0.3, 550.0), (-0.3, 549.87), (-6.566, 549.84335), (-6.566, 549.71335), (-0.3, 549.87), (-0.3, 549.75), (-6.566, 549.71335), (-6.566, 549.59335)]'
$ws.Range("B552").Value = '[(-0.3, 551.0), (-0.3, 550.87), (-6.624, 550.8419), (-6.624, 550.7119), (-0.3, 550.87), (-0.3, 550.75), (-6.624, 550.7119), (-6.624, 550.5919)]'
$ws.Range("B553").Value = '[(-0.3, 552.0), (-0.3, 551.87), (-6.746, 551.83885), (-6.746, 551.70885), (-0.3, 551.87), (-0.3, 551.75), (-6.746, 551.70885), (-6.746, 551.58885)]'
$ws.Range("B554").Value = '[(-0.3, 553.0), (-0.3, 552.87), (-6.891, 552.835225), (-6.891, 552.705225), (-0.3, 552.87), (-0.3, 552.75), (-6.891, 552.705225), (-6.891, 552.585225)]'
$ws.Range("B555").Value = '[(-0.3, 554.0), (-0.3, 553.87), (-7.062, 553.83095), (-7.062, 553.70095), (-0.3, 553.87), (-0.3, 553.75), (-7.062, 553.70095), (-7.062, 553.58095)]'
$ws.Range("B556").Value = '[(-0.3, 555.0), (-0.3, 554.87), (-10.218, 554.75205), (-10.218, 554.6220500000001), (-0.3, 554.87), (-0.3, 554.75), (-10.218, 554.6220500000001), (-10.218, 554.50205)]'
$ws.Range("B557").Value = '[(-0.3, 556.0), (-0.3, 555.87), (-9.355, 555.773625), (-9.355, 555.643625), (-0.3, 555.87), (-0.3, 555.75), (-9.355, 555.643625), (-9.355, 555.523625)]'
$ws.Range("B558").Value = '[(-0.3, 557.0), (-0.3, 556.87), (-8.685, 556.790375), (-8.685, 556.660375), (-0.3, 556.87), (-0.3, 556.75), (-8.685, 556.660375), (-8.685, 556.540375)]'
$ws.Range("B559").Value = '[(-0.3, 558.0), (-0.3, 557.87), (-8.126, 557.80435), (-8.126, 557.67435), (-0.3, 557.87), (-0.3, 557.75), (-8.126, 557.67435), (-8.126, 557.55435)]'
$ws.Range("B560").Value = '[(-0.3, 559.0), (-0.3, 558.87), (-7.618, 558.81705), (-7.618, 558.68705), (-0.3, 558.87), (-0.3, 558.75), (-7.618, 558.68705), (-7.618, 558.56705)]'
$ws.Range("B561").Value = '[0]'
$ws.Range("B562").Value = '[(-0.3, 561.0), (-0.3, 560.81), (-7.039, 560.831525), (-7.039, 560.641525), (-0.3, 560.81), (-0.3, 560.63), (-7.039, 560.641525), (-7.039, 560.511525)]'
$ws.Range("B563").Value = '[(-0.3, 562.0), (-0.3, 561.74), (-6.683, 561.840425), (-6.683, 561.580425), (-0.3, 561.74), (-0.3, 561.48), (-6.683, 561.580425), (-6.683, 561.420425)]'
$ws.Range("B564").Value = '[(-0.3, 563.0), (-0.3, 562.68), (-6.372, 562.852504150004), (-6.372, 562.612504150004), (-0.3, 562.68), (-0.3, 562.36), (-6.372, 562.612504150004), (-6.372, 562.372504150004)]'
$ws.Range("B565").Value = '[(-0.3, 564.0), (-0.3, 563.79), (-6.054, 564.0040787350335), (-6.054, 563.7940787350335), (-0.3, 563.79), (-0.3, 563.58), (-6.054, 563.7940787350335), (-6.054, 563.5840787350336), (-0.3, 563.58), (-0.3, 563.38), (-6.054, 563.5840787350336), (-6.054, 563.4340787350335)]'
$ws.Range("B566").Value = '[(-0.3, 565.0), (-0.3, 564.84), (-5.826, 565.1394990070415), (-5.826, 564.9794990070416), (-0.3, 564.84), (-0.3, 564.68), (-5.826, 564.9794990070416), (-5.826, 564.8194990070415), (-0.3, 564.68), (-0.3, 564.52), (-5.826, 564.8194990070415), (-5.826, 564.6594990070415), (-0.3, 564.52), (-0.3, 564.38), (-5.826, 564.6594990070415), (-5.826, 564.5294990070415)]'
$ws.Range("B567").Value = '[(-0.3, 566.0), (-0.3, 565.8), (-5.72, 566.1834876710203), (-5.72, 565.9834876710203), (-0.3, 565.8), (-0.3, 565.6), (-5.72, 565.9834876710203), (-5.72, 565.7834876710203), (-0.3, 565.6), (-0.3, 565.46), (-5.72, 565.7834876710203), (-5.72, 565.6034876710203)]'
$ws.Range("B568").Value = '[(-0.3, 567.0), (-0.3, 566.8), (-5.575, 567.2239950031015), (-5.575, 567.0239950031015), (-0.3, 566.8), (-0.3, 566.6), (-5.575, 567.0239950031015), (-5.575, 566.8239950031016), (-0.3, 566.6), (-0.3, 566.47), (-5.575, 566.8239950031016), (-5.575, 566.6239950031015)]'
$ws.Range("B569").Value = '[(-0.3, 568.0), (-0.3, 567.8), (-5.352, 568.2580218142435), (-5.352, 568.0580218142435), (-0.3, 567.8), (-0.3, 567.6), (-5.352, 568.0580218142435), (-5.352, 567.8580218142436), (-0.3, 567.6), (-0.3, 567.49), (-5.352, 567.8580218142436), (-5.352, 567.6680218142435)]'
$ws.Range("B570").Value = '[(-0.3, 569.0), (-0.3, 568.7), (-5.076, 569.2850455084694), (-5.076, 568.9850455084694), (-0.3, 568.7), (-0.3, 568.53), (-5.076, 568.9850455084694), (-5.076, 568.6850455084693)]'
$ws.Range("B571").Value = '[(-0.3, 570.0), (-0.3, 569.78), (-4.808, 570.3078629924009), (-4.808, 569.997862992401), (-0.3, 569.78), (-0.3, 569.56), (-4.808, 569.997862992401), (-4.808, 569.6878629924009)]'
$ws.Range("B572").Value = '[(-0.3, 571.0), (-0.3, 570.78), (-4.475, 571.29225), (-4.475, 570.99225), (-0.3, 570.78), (-0.3, 570.56), (-4.475, 570.99225), (-4.475, 570.69225)]'
$ws.Range("B573").Value = '[(-0.3, 572.0), (-0.3, 571.76), (-4.178, 572.27146), (-4.178, 571.95146), (-0.3, 571.76), (-0.3, 571.52), (-4.178, 571.95146), (-4.178, 571.6314600000001)]'
$ws.Range("B574").Value = '[(-0.3, 573.0), (-0.3, 572.72), (-3.883, 573.25081), (-3.883, 572.97081), (-0.3, 572.72), (-0.3, 572.55), (-3.883, 572.97081), (-3.883, 572.69081)]'
$ws.Range("B575").Value = '[(-0.3, 574.0), (-0.3, 573.75), (-3.535, 574.22645), (-3.535, 573.97645), (-0.3, 573.75), (-0.3, 573.58), (-3.535, 573.97645), (-3.535, 573.72645)]'
$ws.Range("B576").Value = '[(-0.3, 575.0), (-0.3, 574.79), (-3.088, 575.19516), (-3.088, 574.98516), (-0.3, 574.79), (-0.3, 574.64), (-3.088, 574.98516), (-3.088, 574.77516)]'
$ws.Range("B577").Value = '[(-0.3, 576.0), (-0.3, 575.77), (-2.856, 576.17892), (-2.856, 575.9489199999999), (-0.3, 575.77), (-0.3, 575.64), (-2.856, 575.9489199999999), (-2.856, 575.7289199999999)]'
$ws.Range("B578").Value = '[(-0.3, 577.0), (-0.3, 576.76), (-2.574, 577.15918), (-2.574, 576.91918), (-0.3, 576.76), (-0.3, 576.58), (-2.574, 576.91918), (-2.574, 576.67918)]'
$ws.Range("B605").Value = '[(-0.3, 604.0), (-0.3, 603.76), (-5.788, 603.8628), (-5.788, 603.7428), (-0.3, 603.76), (-0.3, 603.52), (-5.788, 603.7428), (-5.788, 603.6228)]'
$ws.Range("B622").Value = '[(-0.3, 621.0), (-0.3, 620.86), (-5.488, 620.8210012365475), (-5.488, 620.6810012365476), (-0.3, 620.86), (-0.3, 620.74), (-5.488, 620.6810012365476), (-5.488, 620.5710012365475)]'
$ws.Range("B635").Value = '[(-0.3, 634.0), (-0.3, 633.87), (-7.796, 633.73764), (-7.796, 633.6076400000001), (-0.3, 633.87), (-0.3, 633.75), (-7.796, 633.6076400000001), (-7.796, 633.48764)]'
$ws.Range("B636").Value = '[(-0.3, 635.0), (-0.3, 634.87), (-6.877, 634.769805), (-6.877, 634.639805), (-0.3, 634.87), (-0.3, 634.75), (-6.877, 634.639805), (-6.877, 634.519805)]'
$ws.Range("B637").Value = '[(-0.3, 636.0), (-0.3, 635.87), (-6.963, 635.766795), (-6.963, 635.636795), (-0.3, 635.87), (-0.3, 635.75), (-6.963, 635.636795), (-6.963, 635.516795)]'
$ws.Range("B638").Value = '[(-0.3, 637.0), (-0.3, 636.87), (-7.004, 636.76536), (-7.004, 636.63536), (-0.3, 636.87), (-0.3, 636.75), (-7.004, 636.63536), (-7.004, 636.51536)]'
$ws.Range("B639").Value = '[(-0.3, 638.0), (-0.3, 637.87), (-7.037, 637.764205), (-7.037, 637.634205), (-0.3, 637.87), (-0.3, 637.75), (-7.037, 637.634205), (-7.037, 637.514205)]'
$ws.Range("B640").Value = '[(-0.3, 639.0), (-0.3, 638.87), (-7.102, 638.76193), (-7.102, 638.63193), (-0.3, 638.87), (-0.3, 638.75), (-7.102, 638.63193), (-7.102, 638.51193)]'
$ws.Range("B641").Value = '[(-0.3, 640.0), (-0.3, 639.87), (-7.456, 639.77809244), (-7.456, 639.64809244), (-0.3, 639.87), (-0.3, 639.75), (-7.456, 639.64809244), (-7.456, 639.52809244)]'
$ws.Range("B642").Value = '[(-0.3, 641.0), (-0.3, 640.87), (-7.614, 640.7637578), (-7.614, 640.6337578), (-0.3, 640.87), (-0.3, 640.75), (-7.614, 640.6337578), (-7.614, 640.5137578)]'
$ws.Range("B643").Value = '[(-0.3, 641.6355600000002), (-0.3, 641.5055600000002), (-7.624, 641.4084427600002), (-7.624, 641.2784427600002), (-0.3, 641.5055600000002), (-0.3, 641.3855600000002), (-7.624, 641.2784427600002), (-7.624, 641.1584427600002)]'
